# =========================================================================
# Add the 'Duke 18-19' worksheet (COUNTER JR2 report for Duke University
# Press titles on the Silverchair platform) as the 4th / last tab, and make
# it the active sheet -- mirrors the existing 'PQ 18-19' report tab.
# =========================================================================
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Duke 18-19"

# --- Column widths ---
$ws.Columns("A:B").ColumnWidth = 35.166666666666664   # -> ~36 chars
$ws.Columns("C:J").ColumnWidth = 9.833333333333334    # -> ~10.71 chars
$ws.Columns("K:U").ColumnWidth = 4.666666666666667    # -> ~5.43 chars

# --- Row heights ---
$ws.Rows(1).RowHeight = 15.95
$ws.Rows(2).RowHeight = 15.95
$ws.Rows(3).RowHeight = 15.95
$ws.Rows(4).RowHeight = 15.95
$ws.Rows(5).RowHeight = 15.95
$ws.Rows(6).RowHeight = 15.95
$ws.Rows(7).RowHeight = 12
$ws.Rows(8).RowHeight = 45
$ws.Rows(9).RowHeight = 15
$ws.Rows(10).RowHeight = 15
$ws.Rows(11).RowHeight = 15
$ws.Rows(12).RowHeight = 15
$ws.Rows(13).RowHeight = 15
$ws.Rows(14).RowHeight = 15
$ws.Rows(15).RowHeight = 15
$ws.Rows(16).RowHeight = 15
$ws.Rows(17).RowHeight = 15
$ws.Rows(18).RowHeight = 15
$ws.Rows(19).RowHeight = 15
$ws.Rows(20).RowHeight = 15

# --- Cell formatting (font / fill / alignment) per region ---

# Style group 3: title block (rows 1-7, col A, plus B1) -- bold 9pt Arial,
# left/top aligned, no fill
$rngTitle = $ws.Range("A1:A7,B1")
$rngTitle.Font.Name = "Arial"
$rngTitle.Font.Size = 9
$rngTitle.Font.Bold = $true
$rngTitle.HorizontalAlignment = -4131
$rngTitle.VerticalAlignment = -4160

# Style group 5: header row (row 8), non-wrapping columns -- bold 9pt Arial,
# dark-blue fill, left aligned
$rngHdr1 = $ws.Range("A8:D8,F8:J8")
$rngHdr1.Font.Name = "Arial"
$rngHdr1.Font.Size = 9
$rngHdr1.Font.Bold = $true
$rngHdr1.Interior.Color = 11439959
$rngHdr1.HorizontalAlignment = -4131

# Style group 6: header row (row 8), wrapping columns -- bold 9pt Arial,
# dark-blue fill, left aligned, wrap text
$rngHdr2 = $ws.Range("E8,K8:U8")
$rngHdr2.Font.Name = "Arial"
$rngHdr2.Font.Size = 9
$rngHdr2.Font.Bold = $true
$rngHdr2.Interior.Color = 11439959
$rngHdr2.HorizontalAlignment = -4131
$rngHdr2.WrapText = $true

# Style group 7: 'Total for all journals' rows (9-10) -- regular 9pt Arial,
# light-blue fill, left/top aligned, wrap text
$rngTotal = $ws.Range("A9:U10")
$rngTotal.Font.Name = "Arial"
$rngTotal.Font.Size = 9
$rngTotal.Interior.Color = 14602940
$rngTotal.HorizontalAlignment = -4131
$rngTotal.VerticalAlignment = -4160
$rngTotal.WrapText = $true

# Style group 8: data rows (11-20) -- regular 9pt Arial, no fill,
# left/top aligned, wrap text
$rngData = $ws.Range("A11:U20")
$rngData.Font.Name = "Arial"
$rngData.Font.Size = 9
$rngData.HorizontalAlignment = -4131
$rngData.VerticalAlignment = -4160
$rngData.WrapText = $true

# --- Cell values ---
$ws.Range("A1").Value = 'Journal Report 2 (R4)'
$ws.Range("B1").Value = 'Access Denied to  Full-Text Article by Month, Journal and Category'

$ws.Range("A2").Value = 'FLORIDA STATE UNIV'

$ws.Range("A3").Value = ' '

$ws.Range("A4").Value = 'Period covered by Report:'

$ws.Range("A5").Value = '2018-07-01 to 2019-06-30'

$ws.Range("A6").Value = 'Date run:'

$ws.Range("A7").Value = '2019-09-05'

$ws.Range("A8").Value = 'Journal'
$ws.Range("B8").Value = 'Publisher'
$ws.Range("C8").Value = 'Platform'
$ws.Range("D8").Value = 'Journal DOI'
$ws.Range("E8").Value = 'Proprietary Identifier'
$ws.Range("F8").Value = 'Print ISSN'
$ws.Range("G8").Value = 'Online ISSN'
$ws.Range("H8").Value = 'Access Denied Category'
$ws.Range("I8").Value = 'Reporting Period Total'
$ws.Range("J8").Value = 'Jul-2018'
$ws.Range("K8").Value = 'Aug-2018'
$ws.Range("L8").Value = 'Sep-2018'
$ws.Range("M8").Value = 'Oct-2018'
$ws.Range("N8").Value = 'Nov-2018'
$ws.Range("O8").Value = 'Dec-2018'
$ws.Range("P8").Value = 'Jan-2019'
$ws.Range("Q8").Value = 'Feb-2019'
$ws.Range("R8").Value = 'Mar-2019'
$ws.Range("S8").Value = 'Apr-2019'
$ws.Range("T8").Value = 'May-2019'
$ws.Range("U8").Value = 'Jun-2019'

$ws.Range("A9").Value = 'Total for all journals'
$ws.Range("C9").Value = 'Silverchair'
$ws.Range("H9").Value = 'Access Denied: concurrent/simultaneous user license limit exceeded'

$ws.Range("A10").Value = 'Total for all journals'
$ws.Range("C10").Value = 'Silverchair'
$ws.Range("H10").Value = 'Access Denied: content item not licensed'

$ws.Range("A11").Value = 'Camera Obscura: Feminism, Culture, and Media Studies'
$ws.Range("B11").Value = 'Duke University Press'
$ws.Range("C11").Value = 'Silverchair'
$ws.Range("F11").Value = '0270-5346'
$ws.Range("G11").Value = '1529-1510'
$ws.Range("H11").Value = 'Access Denied: concurrent/simultaneous user license limit exceeded'
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 0
$ws.Range("T11").Value = 0
$ws.Range("U11").Value = 0

$ws.Range("A12").Value = 'Camera Obscura: Feminism, Culture, and Media Studies'
$ws.Range("B12").Value = 'Duke University Press'
$ws.Range("C12").Value = 'Silverchair'
$ws.Range("F12").Value = '0270-5346'
$ws.Range("G12").Value = '1529-1510'
$ws.Range("H12").Value = 'Access Denied: content item not licensed'
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0
$ws.Range("R12").Value = 0
$ws.Range("S12").Value = 0
$ws.Range("T12").Value = 0
$ws.Range("U12").Value = 0

$ws.Range("A13").Value = 'Environmental Humanities'
$ws.Range("B13").Value = 'Duke University Press'
$ws.Range("C13").Value = 'Silverchair'
$ws.Range("F13").Value = '2201-1919'
$ws.Range("G13").Value = '2201-1919'
$ws.Range("H13").Value = 'Access Denied: concurrent/simultaneous user license limit exceeded'
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0
$ws.Range("R13").Value = 0
$ws.Range("S13").Value = 0
$ws.Range("T13").Value = 0
$ws.Range("U13").Value = 0

$ws.Range("A14").Value = 'Environmental Humanities'
$ws.Range("B14").Value = 'Duke University Press'
$ws.Range("C14").Value = 'Silverchair'
$ws.Range("F14").Value = '2201-1919'
$ws.Range("G14").Value = '2201-1919'
$ws.Range("H14").Value = 'Access Denied: content item not licensed'
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0
$ws.Range("R14").Value = 0
$ws.Range("S14").Value = 0
$ws.Range("T14").Value = 0
$ws.Range("U14").Value = 0

$ws.Range("A15").Value = 'Journal of Health Politics, Policy and Law'
$ws.Range("B15").Value = 'Duke University Press'
$ws.Range("C15").Value = 'Silverchair'
$ws.Range("F15").Value = '0361-6878'
$ws.Range("G15").Value = '1527-1927'
$ws.Range("H15").Value = 'Access Denied: concurrent/simultaneous user license limit exceeded'
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0
$ws.Range("R15").Value = 0
$ws.Range("S15").Value = 0
$ws.Range("T15").Value = 0
$ws.Range("U15").Value = 0

$ws.Range("A16").Value = 'Journal of Health Politics, Policy and Law'
$ws.Range("B16").Value = 'Duke University Press'
$ws.Range("C16").Value = 'Silverchair'
$ws.Range("F16").Value = '0361-6878'
$ws.Range("G16").Value = '1527-1927'
$ws.Range("H16").Value = 'Access Denied: content item not licensed'
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 0
$ws.Range("R16").Value = 0
$ws.Range("S16").Value = 0
$ws.Range("T16").Value = 0
$ws.Range("U16").Value = 0

$ws.Range("A17").Value = 'Journal of Korean Studies'
$ws.Range("B17").Value = 'Duke University Press'
$ws.Range("C17").Value = 'Silverchair'
$ws.Range("F17").Value = '2158-1665'
$ws.Range("G17").Value = '0731-1613'
$ws.Range("H17").Value = 'Access Denied: concurrent/simultaneous user license limit exceeded'
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 0
$ws.Range("R17").Value = 0
$ws.Range("S17").Value = 0
$ws.Range("T17").Value = 0
$ws.Range("U17").Value = 0

$ws.Range("A18").Value = 'Journal of Korean Studies'
$ws.Range("B18").Value = 'Duke University Press'
$ws.Range("C18").Value = 'Silverchair'
$ws.Range("F18").Value = '2158-1665'
$ws.Range("G18").Value = '0731-1613'
$ws.Range("H18").Value = 'Access Denied: content item not licensed'
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 0
$ws.Range("R18").Value = 0
$ws.Range("S18").Value = 0
$ws.Range("T18").Value = 0
$ws.Range("U18").Value = 0

$ws.Range("A19").Value = 'the minnesota review'
$ws.Range("B19").Value = 'Duke University Press'
$ws.Range("C19").Value = 'Silverchair'
$ws.Range("F19").Value = '0026-5667'
$ws.Range("G19").Value = '2157-4189'
$ws.Range("H19").Value = 'Access Denied: concurrent/simultaneous user license limit exceeded'
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 0
$ws.Range("R19").Value = 0
$ws.Range("S19").Value = 0
$ws.Range("T19").Value = 0
$ws.Range("U19").Value = 0

$ws.Range("A20").Value = 'the minnesota review'
$ws.Range("B20").Value = 'Duke University Press'
$ws.Range("C20").Value = 'Silverchair'
$ws.Range("F20").Value = '0026-5667'
$ws.Range("G20").Value = '2157-4189'
$ws.Range("H20").Value = 'Access Denied: content item not licensed'
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 0
$ws.Range("R20").Value = 0
$ws.Range("S20").Value = 0
$ws.Range("T20").Value = 0
$ws.Range("U20").Value = 0

# --- Page setup ---
$ws.PageSetup.Orientation = 2
$ws.PageSetup.PrintGridlines = $true

# --- View: make the new sheet the active / selected tab, matching the
# original author's last on-screen position ---
$ws.Range("B34").Select()
$ws.Activate()

